$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column H ("property_category") between "total" (G) and "date" (old H),
# pushing date / legislator_name / legislator_id one column to the right.
$ws.Columns("H:H").Insert()

# Header row
$ws.Range("H1").Value = "property_category"

# Data rows - every stock entry on this sheet is a "stock" category
$ws.Range("H2").Value = "stock"
$ws.Range("H3").Value = "stock"
$ws.Range("H4").Value = "stock"
